# Update Tutorial 6 solution for 2001ME38: change date separators from "/"
# to "-" and correct the Total/Real/Invalid/Absent attendance counts for
# the corresponding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, Date(with dashes), D(Total Attendance Count), E(Real), F(Duplicate), G(Invalid), H(Absent)
$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; F = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 6;  Date = "08-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 11; Date = "25-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 12; Date = "29-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 13; Date = "01-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 14; Date = "05-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 15; Date = "08-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 16; Date = "12-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $cell = $ws.Range("A$i")

    # Some of the new date strings (e.g. "01-08-2022") are ambiguous and
    # would otherwise be auto-recognized by Excel as a date value (since
    # the day component is <= 12). Force a Text number format while we
    # assign the value so the literal dash-separated string is preserved,
    # then restore the default "Normal" style so the cell formatting
    # matches the rest of the column.
    $cell.NumberFormat = "@"
    $cell.Value = $r.Date
    $cell.Style = "Normal"

    $ws.Range("D$i").Value = $r.D
    $ws.Range("E$i").Value = $r.E
    $ws.Range("F$i").Value = $r.F
    $ws.Range("G$i").Value = $r.G
    $ws.Range("H$i").Value = $r.H
}
